# Generate Report for Handback
# Updates the "Latest HO / Correspond Handoff / Handback" datetime values
# for the second data row (file aa8163b8-...) in each report sheet, which
# previously were (incorrectly) duplicated from the first data row.

$wb = $excel.ActiveWorkbook

# --- Overview sheet -------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G3").Value = "2016-08-25 04:46:45"

# --- zh-cn sheet ------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H3").Value = "2016-08-25 04:46:40"
$zhcn.Range("K3").Value = "2016-08-25 04:46:57"

# --- de-de sheet ------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("H3").Value = "2016-08-25 04:46:45"
$dede.Range("K3").Value = "2016-08-25 04:47:09"
